# The workbook has a single worksheet ("Sheet1"); grab it explicitly rather
# than relying solely on ActiveSheet, then append the new data row (row 81)
# at the bottom of the existing table (which ran through row 80).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A81").Value = "2024-10-22 00:00:00"
$ws.Range("B81").Value = 73350
$ws.Range("C81").Value = 10276.71
$ws.Range("D81").Value = 9094.43
$ws.Range("E81").Value = 7.1208
